$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values (IMG updates + lvl updates)
$ws.Range("BA1").Value = 3
$ws.Range("BB1").Value = 3
$ws.Range("BC1").Value = 3
$ws.Range("CY1").Value = 0
$ws.Range("DB1").Value = 6

$ws.Range("BA2").Value = 3
$ws.Range("BB2").Value = 8
$ws.Range("BC2").Value = 8
$ws.Range("CY2").Value = 0
$ws.Range("DB2").Value = 4
$ws.Range("DF2").Value = 4
$ws.Range("DG2").Value = 4
$ws.Range("DH2").Value = 4
$ws.Range("DI2").Value = 4
$ws.Range("DJ2").Value = 4
$ws.Range("DK2").Value = 4
$ws.Range("DL2").Value = 4

$ws.Range("BA3").Value = 8
$ws.Range("BD3").Value = 0

$ws.Range("AZ6").Value = 3
$ws.Range("BA6").Value = 0

$ws.Range("AZ7").Value = 3

$ws.Range("DE8").Value = 7
$ws.Range("DF8").Value = 0
$ws.Range("DK8").Value = 8
$ws.Range("DL8").Value = 8

$ws.Range("DM9").Value = 1
$ws.Range("DN9").Value = 1
$ws.Range("DO9").Value = 1
$ws.Range("DP9").Value = 1
$ws.Range("DQ9").Value = 1
$ws.Range("DR9").Value = 1
$ws.Range("DS9").Value = 1

# Update the view's top-left cell and selection to match the edited area
$null = $ws.Range("AO1").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 41
$null = $ws.Range("BB6").Select()
